# edit.ps1 -- PowerPoint COM-interop script (PowerShell-style) that
# reproduces the authored change:
#   - bump the cached "last saved" date field text on the slide
#     master + every slide layout from 26-11-2022 to 30-11-2022
#   - append a new "2. Environment Setup" slide (Title and Content
#     layout) listing the environment-setup tools
#
# $ppt.ActivePresentation is already open as $p.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Refresh the cached date-field text ("26-11-2022" -> "30-11-2022")
#    on the slide master and on every slide layout's Date Placeholder.
# ---------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "26-11-2022") {
                    $tr.LanguageID = "en-IN"
                    $tr.Text = "30-11-2022"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------
# 2) Add the new "2. Environment Setup" slide at the end (index 7)
#    using the "Title and Content" layout (ppLayoutText = 2), which
#    matches the layout used by the other slides in the deck.
# ---------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.LanguageID = "en-IN"
$titleShape.TextFrame.TextRange.Text = "2. Environment Setup"

$contentShape = $newSlide.Shapes.Item(2)
$contentShape.TextFrame.TextRange.LanguageID = "en-IN"
$contentShape.TextFrame.TextRange.Text = "Node`rNpm`rJava`rAndroid Studio`rExpo`rVs Code"
